$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text-number-format on the cells we touch so values like
# "27.20" or "0.002184" keep their exact text representation
# (trailing zeros, percent signs) instead of being auto-coerced to
# floating point numbers by the Value setter's type inference.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "261.53"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "1.14%"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "27.20"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "1.20%"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "4.703"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "0.57%"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.06195"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "3.28%"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "6.716"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "0.76%"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "-0.76%"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.9135"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "-1.26%"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.1413"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "1.58%"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.04537"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "-0.18%"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07087"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "0.75%"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.03141"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "0.74%"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.09039"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "-1.00%"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.001530"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "0.41%"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.0006166"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "1.94%"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.006043"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "0.28%"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.163"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "-0.07%"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "1.32%"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.3106"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "-0.18%"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "1.02%"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.095"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "-1.07%"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.001217"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "0.10%"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "-5.72%"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "0.06%"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.03936"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "2.48%"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.1114"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "0.01%"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.004124"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "7.39%"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.002184"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "-9.73%"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.01382"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "-8.50%"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00005150"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "0.71%"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "0.07%"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "0.07%"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "0.07%"
